$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Unprotect()

# Update the confidentiality/date note in A16
$ws.Range("A16").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-04-26 for illustrative purposes only and are subject to change."

# Update Weight (D) and Percent Change (E) values for rows 2-12
$ws.Range("D2").Value = 0.03067474672912161
$ws.Range("E2").Value = 0.003765833618623704

$ws.Range("D3").Value = 0.02390201695170641
$ws.Range("E3").Value = -0.00649651972157772

$ws.Range("D4").Value = 0.05111423375315943
$ws.Range("E4").Value = -0.01059571462208608

$ws.Range("D5").Value = 0.1389565354264186
$ws.Range("E5").Value = 0.00128390306531867

$ws.Range("D6").Value = 0.02833127060394097
$ws.Range("E6").Value = 0.007009345794392496

$ws.Range("D7").Value = 0.1209756865384459
$ws.Range("E7").Value = 0.004211569428960615

$ws.Range("D8").Value = 0.1007096504939164
$ws.Range("E8").Value = -0.001851166234727808

$ws.Range("D9").Value = 0.02791003251084924
$ws.Range("E9").Value = 0.007115135834411257

$ws.Range("D10").Value = 0.1202121482470183
$ws.Range("E10").Value = 0.003729878288182276

$ws.Range("D11").Value = 0.2549469812272943
$ws.Range("E11").Value = 0.008265916285613795

$ws.Range("D12").Value = 0.1022666975181289
$ws.Range("E12").Value = 0.003101376235704612

# Row 13 (Total) - only E13 changes, D13 stays 1
$ws.Range("E13").Value = 0.00319020009045401

# Restore sheet protection (the sheet was protected before the edit)
$ws.Protect()
